# Trade #42 closed at 2026-02-18 00:15:59 - unknown UNKNOWN +0.000%
#
# This updates the live trading results workbook:
#  - Summary sheet roll-up numbers (Total P&L %, Total Trades, Win Rate %)
#  - Strategy Status row for "momentum" (Trades, Win Rate %)
#  - "All Trades" sheet: closes the existing OPEN momentum trade (#70, row 71)
#    and appends a brand-new OPEN momentum trade (#99)
#  - "momentum" sheet: same two updates, mirrored with its own column layout

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Helper: write a text value into a cell while stopping Excel from quietly
# re-interpreting date/time-shaped strings ("2026-02-18", "00:15:53") as
# serial date numbers. We force a Text number format for the write, then
# strip the formatting back off so the cell ends up as a plain, unstyled
# text cell again (matching every other text cell in this sheet).
# ---------------------------------------------------------------------------
function Set-TextValue {
    param($cell, [string]$value)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.ClearFormats()
}

# ---------------------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B5").Value = 0.16     # Total P&L %
$wsSummary.Range("B6").Value = 70       # Total Trades
$wsSummary.Range("B9").Value = 51.43    # Win Rate %

# ---------------------------------------------------------------------------
# Strategy Status sheet - "momentum" row (row 11)
# ---------------------------------------------------------------------------
$wsStatus = $wb.Worksheets.Item("Strategy Status")
$wsStatus.Range("D11").Value = 7        # Trades
$wsStatus.Range("G11").Value = 14.29    # Win Rate %

# ---------------------------------------------------------------------------
# All Trades sheet
#   columns: A Trade#, B Date, C Time, D Strategy, E Side, F Entry Price,
#            G Exit Price, H Status, I P&L %, J P&L $, K Capital After,
#            L Exit Reason, M Duration (min), N Entry Slippage (bps),
#            O Exit Slippage (bps), P Confidence, Q Entry Reason
# ---------------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("All Trades")

# Close trade #70 (row 71): was OPEN with no exit info, now CLOSED.
$wsAll.Range("G71").Value = 0.99
$wsAll.Range("H71").Value = "CLOSED"
$wsAll.Range("K71").Value = 99.68000000000001
Set-TextValue $wsAll.Range("L71") "early_exit"
$wsAll.Range("M71").Value = 0.12

# Append new trade #99 as row 100 (still OPEN).
Set-TextValue $wsAll.Range("B100") "2026-02-18"
Set-TextValue $wsAll.Range("C100") "00:15:53"
$wsAll.Range("A100").Value = 99
$wsAll.Range("D100").Value = "momentum"
$wsAll.Range("E100").Value = "UP"
$wsAll.Range("F100").Value = 0.99
$wsAll.Range("H100").Value = "OPEN"
$wsAll.Range("I100").Value = 0
$wsAll.Range("J100").Value = 0
$wsAll.Range("K100").Value = 99.6787371310913
$wsAll.Range("M100").Value = 0
$wsAll.Range("N100").Value = 0
$wsAll.Range("O100").Value = 0
$wsAll.Range("P100").Value = 0.9
$wsAll.Range("Q100").Value = "Upward momentum: 71.186% over 10 samples"

# ---------------------------------------------------------------------------
# "momentum" sheet (strategy-specific view, different column order):
#   columns: A Trade#, B Date, C Time, D Strategy, E Side, F Entry Price,
#            G Exit Price, H Status, I P&L %, J P&L $, K Capital After,
#            L Entry Slippage (bps), M Exit Slippage (bps), N Confidence,
#            O Entry Reason, P Exit Reason, Q Duration (min)
# ---------------------------------------------------------------------------
$wsMom = $wb.Worksheets.Item("momentum")

# Close trade #70 (row 8): was OPEN with no exit info, now CLOSED.
$wsMom.Range("G8").Value = 0.99
$wsMom.Range("H8").Value = "CLOSED"
$wsMom.Range("K8").Value = 99.68000000000001
Set-TextValue $wsMom.Range("P8") "early_exit"
$wsMom.Range("Q8").Value = 0.12

# Append new trade #99 as row 23 (still OPEN).
Set-TextValue $wsMom.Range("B23") "2026-02-18"
Set-TextValue $wsMom.Range("C23") "00:15:53"
$wsMom.Range("A23").Value = 99
$wsMom.Range("D23").Value = "momentum"
$wsMom.Range("E23").Value = "UP"
$wsMom.Range("F23").Value = 0.99
$wsMom.Range("H23").Value = "OPEN"
$wsMom.Range("I23").Value = 0
$wsMom.Range("J23").Value = 0
$wsMom.Range("K23").Value = 99.6787371310913
$wsMom.Range("L23").Value = 0
$wsMom.Range("M23").Value = 0
$wsMom.Range("N23").Value = 0.9
$wsMom.Range("O23").Value = "Upward momentum: 71.186% over 10 samples"
$wsMom.Range("Q23").Value = 0
